$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: insert a brand-new plain ("Normal") paragraph right after the
# paragraph object $anchor, put $text in it, and return the new
# paragraph object so callers can keep chaining.
# ---------------------------------------------------------------------
function Add-PlainParagraphAfter($anchor, $text) {
    $anchor.Range.InsertParagraphAfter()
    $newPara = $anchor.Next()
    # Only force the "Normal" style when the anchor itself is not already
    # plain text - a paragraph inserted after an existing "Normal"
    # paragraph already serialises without any pPr/pStyle at all, which
    # is the cleanest match for a genuinely un-styled paragraph.
    if ($anchor.Style.NameLocal -ne "Normal") {
        $newPara.Style = "Normal"
    }
    if ($text -ne $null -and $text -ne "") {
        $newPara.Range.Text = $text
    }
    return $newPara
}

# Helper: insert a brand-new Heading-1 paragraph right after $anchor.
function Add-Heading1After($anchor, $text) {
    $anchor.Range.InsertParagraphAfter()
    $newPara = $anchor.Next()
    $newPara.Style = "Heading 1"
    $newPara.Range.Text = $text
    return $newPara
}

# ---------------------------------------------------------------------
# Step C - drop the leftover boilerplate at the tail of the document:
# the extra "Heading, level 1" heading, the intense-quote paragraph,
# both sample list items and the sample Qty/Id/Desc table.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Heading, level 1") | Out-Null
$p = $rng.Paragraphs(1)
$p.Range.Delete()

$rng = $d.Content
$rng.Find.Execute("Intense quote") | Out-Null
$p = $rng.Paragraphs(1)
$p.Range.Delete()

$rng = $d.Content
$rng.Find.Execute("first item in unordered list") | Out-Null
$p = $rng.Paragraphs(1)
$p.Range.Delete()

$rng = $d.Content
$rng.Find.Execute("first item in ordered list") | Out-Null
$p = $rng.Paragraphs(1)
$p.Range.Delete()

if ($d.Tables.Count -gt 0) {
    $d.Tables(1).Delete()
}

# ---------------------------------------------------------------------
# Step B - the demo "bold / italic" paragraph becomes a single plain
# run with the contract amount in words.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("A plain paragraph having some") | Out-Null
$p = $rng.Paragraphs(1)
$contentRange = $d.Range($p.Range.Start, $p.Range.End - 1)
$contentRange.Text = "сто тысяч, пятьсот рублей"

# ---------------------------------------------------------------------
# Step A - walk the numbered headings top to bottom, fix their text and
# hang the newly-authored detail paragraphs off each one.
# ---------------------------------------------------------------------

# 1. Номер договора
$rng = $d.Content
$rng.Find.Execute("1. Номер договора bred", $true, $false, $false, $false, $false, $true, 1, $false, "1. Номер договора №12320161023170", 2) | Out-Null

# 2. Юридическое лицо со стороны исполнителя
$rng = $d.Content
$rng.Find.Execute("2. Юридическое лицо со стороны исполнителя 2", $true, $false, $false, $false, $false, $true, 1, $false, "2. Юридическое лицо со стороны исполнителя", 2) | Out-Null
$p = $rng.Paragraphs(1)
Add-PlainParagraphAfter $p "ООО Промвад Софт" | Out-Null

# 3. Выбор подписанта со стороны исполнителя
$rng = $d.Content
$rng.Find.Execute("3. Выбор подписанта со стороны исполнителя bred", $true, $false, $false, $false, $false, $true, 1, $false, "3. Выбор подписанта со стороны исполнителя 2", 2) | Out-Null
$p = $rng.Paragraphs(1)
Add-PlainParagraphAfter $p "Ковалев С.Н." | Out-Null

# 4. Место составления договора - text unchanged, nothing to do here.

# The paragraph that used to read "5. Дата договора 1111" loses its
# Heading-1 styling entirely and becomes a plain "Минск" paragraph.
# Build the replacement fresh (after "4. Место составления...") and
# delete the old heading paragraph so no pPr/pStyle survives.
$rng = $d.Content
$rng.Find.Execute("4. Место составления договора 1") | Out-Null
$p4 = $rng.Paragraphs(1)
Add-PlainParagraphAfter $p4 "Минск" | Out-Null

$rng = $d.Content
$rng.Find.Execute("5. Дата договора 1111") | Out-Null
$oldP5 = $rng.Paragraphs(1)
$oldP5.Range.Delete()

# 5. Дата договора (was "6. Информация о Заказчике 123 321 132")
$rng = $d.Content
$rng.Find.Execute("6. Информация о Заказчике 123 321 132", $true, $false, $false, $false, $false, $true, 1, $false, "5. Дата договора 123 20161023170", 2) | Out-Null
$p = $rng.Paragraphs(1)
$p = Add-Heading1After $p "6. Информация о Заказчике:"
$p = Add-PlainParagraphAfter $p "sad"
$p = Add-PlainParagraphAfter $p "ewq"
$p = Add-PlainParagraphAfter $p "qew"

# 7. Информация о представителе Заказчика
$rng = $d.Content
$rng.Find.Execute("7. Информация о представителе Заказчика bred", $true, $false, $false, $false, $false, $true, 1, $false, "7. Информация о представителе Заказчика", 2) | Out-Null
$p = $rng.Paragraphs(1)
$p = Add-PlainParagraphAfter $p "asd"
$p = Add-PlainParagraphAfter $p "dsa"
$p = Add-PlainParagraphAfter $p "sad"
$p = Add-PlainParagraphAfter $p "qwe"
$p = Add-PlainParagraphAfter $p "aad"
$p = Add-PlainParagraphAfter $p "Устава"

# 8. НДС
$rng = $d.Content
$rng.Find.Execute("8. НДС bred", $true, $false, $false, $false, $false, $true, 1, $false, "8. НДС 2", 2) | Out-Null
$p = $rng.Paragraphs(1)
$p = Add-PlainParagraphAfter $p "18% (для РФ)"
$p = Add-PlainParagraphAfter $p ""

# 9. Валюта платежа
$rng = $d.Content
$rng.Find.Execute("9. Валюта платежа bred", $true, $false, $false, $false, $false, $true, 1, $false, "9. Валюта платежа", 2) | Out-Null
$p = $rng.Paragraphs(1)
$p = Add-PlainParagraphAfter $p "Валюта по договору"
$p = Add-PlainParagraphAfter $p "USD"
$p = Add-PlainParagraphAfter $p "Валюта платежа"
$p = Add-PlainParagraphAfter $p "EUR"

# 10. Сумма по договору
$rng = $d.Content
$rng.Find.Execute("10.Сумма по  договору bred", $true, $false, $false, $false, $false, $true, 1, $false, "10.Сумма по  договору {}", 2) | Out-Null
